$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the base value that drives the dependent formulas in row 4
$ws.Range("C4").Value = 1000

# Update the formula for D4 to add 1000 instead of 500
$ws.Range("D4").Formula = "=C4+1000"

# Move the active selection to E19 as recorded in the saved view state
$ws.Range("E19").Select()
